$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark currently sitting after the
#    "Explain Primary Key and Foreign Key ... [2]" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Locate " (many to many)" and find the split point between
#    " (many to m" and "any)" (i.e. right after the "m" in "many").
$search = $d.Content
$search.Find.Execute(" (many to m", $false, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$splitPoint = $search.End

# Re-add the _GoBack bookmark as a zero-length range at that split
# point -- this also breaks the run in two when the document is saved,
# matching Word's behaviour of inserting a bookmark mid-run.
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
